$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 7174
$ws.Range("C3").Value = 155462
$ws.Range("C4").Value = 146576
$ws.Range("C8").Value = 63.62
